{"js": "// The document contained a short C++ code listing (Header.h / source.cpp\n// snippets). The feature described by the listing is finished, so the\n// whole listing is removed, leaving only the single paragraph that holds\n// the `_GoBack` bookmark (its own run/text was already empty).\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Delete every paragraph that still carries text from the old code\n// listing; keep the (already empty) bookmark paragraph untouched so its\n// bookmarkStart/bookmarkEnd survive.\nfor (let i = paragraphs.items.length - 1; i >= 0; i--) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text !== \"\") {\n    paragraph.delete();\n  }\n}\nawait context.sync();\n", "ps1": "# The document contained a short C++ code listing (Header.h / source.cpp\n# snippets). The feature described by the listing is finished, so the\n# whole listing is removed, leaving only the single paragraph that holds\n# the `_GoBack` bookmark (it already had no text of its own).\n$d = $word.ActiveDocument\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    # A paragraph's Range.Text always ends with its paragraph mark (\\r);\n    # an \"empty\" paragraph's text is just that \\r.\n    if ($p.Range.Text -ne \"`r\") {\n        $p.Range.Delete()\n    }\n}\n"}
